$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D6").Value = "2016-01-15 02:59:41"
$wsDe.Range("D6").Value = "2016-01-15 02:59:53"
